# Updates cryptos list cell values per the authored diff.
# Numeric-looking text values (e.g. "565.39") are forced to stay as
# plain text (matching the original inlineStr cells) by prefixing a
# leading apostrophe (Excel's "treat as text" entry marker) and then
# resetting the cell style back to Normal so no stray quotePrefix/
# number-format style is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.487.55'
$ws.Range('D3').Value = '2.440.41'
$ws.Range('E3').Value = '  -1.50%  '
$ws.Range('E4').Value = '  -0.25%  '
$c = $ws.Range('D5')
$c.Value = "'565.39"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.99%  '
$c = $ws.Range('D6')
$c.Value = "'144.23"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -3.37%  '
$c = $ws.Range('D8')
$c.Value = "'0.530"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -2.16%  '
$ws.Range('D9').Value = '2.437.66'
$ws.Range('E9').Value = '  -1.98%  '
$ws.Range('E10').Value = '  -5.60%  '
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('E12').Value = '  -3.09%  '
$ws.Range('E13').Value = '  -2.94%  '
$c = $ws.Range('D14')
$c.Value = "'26.58"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -3.38%  '
$ws.Range('E15').Value = '  -6.13%  '
$ws.Range('D16').Value = '2.875.88'
$ws.Range('E16').Value = '  -2.49%  '
$ws.Range('D17').Value = '62.398.91'
$ws.Range('E17').Value = '  -1.75%  '
$ws.Range('D18').Value = '2.392.54'
$ws.Range('E18').Value = '  -3.86%  '
$ws.Range('E19').Value = '  -4.19%  '
$ws.Range('E20').Value = '  -2.02%  '
$c = $ws.Range('D21')
$c.Value = "'325.04"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('E22').Value = '  -2.63%  '
$c = $ws.Range('D23')
$c.Value = "'2.03"
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +6.26%  '
$ws.Range('E24').Value = '  +0.33%  '
$c = $ws.Range('D25')
$c.Value = "'65.21"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -3.60%  '
$c = $ws.Range('D26')
$c.Value = "'627.16"
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.92%  '
$c = $ws.Range('D27')
$c.Value = "'9.07"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +3.13%  '
$ws.Range('E28').Value = '  -9.15%  '
$ws.Range('E29').Value = '  -1.84%  '
$c = $ws.Range('D30')
$c.Value = "'0.998"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('E32').Value = '  -4.87%  '
$c = $ws.Range('D33')
$c.Value = "'1.86"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -2.55%  '
$c = $ws.Range('D34')
$c.Value = "'0.133"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -7.99%  '
$c = $ws.Range('D35')
$c.Value = "'5.04"
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -3.64%  '
$ws.Range('E37').Value = '  -5.91%  '
$c = $ws.Range('D38')
$c.Value = "'0.376"
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -2.98%  '
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('E40').Value = '  -5.55%  '
$c = $ws.Range('D41')
$c.Value = "'146.06"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('E42').Value = '  -6.75%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D43')
$c.Value = "'42.29"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c = $ws.Range('D44')
$c.Value = "'0.999"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('E45').Value = '  -5.12%  '
$c = $ws.Range('D46')
$c.Value = "'145.65"
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -3.75%  '
$ws.Range('E47').Value = '  -1.93%  '
$c = $ws.Range('D48')
$c.Value = "'20.29"
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -4.65%  '
$c = $ws.Range('D49')
$c.Value = "'0.0526"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -4.78%  '
$c = $ws.Range('D50')
$c.Value = "'0.595"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -2.74%  '
$ws.Range('E51').Value = '  -4.75%  '
